# debug_lines.xlsx - update the "intervals" input column (D) on Sheet1.
# D2, D7 and D8 are the user-entered inputs; column A recalculates
# automatically (A[n] = A[n-1] + D[n-1]) as a consequence.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 1247
$ws.Range("D7").Value = 435
$ws.Range("D8").Value = 315

# leave the selection on D2, matching the saved view state
$ws.Range("D2").Select()
